$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new page link row into the table (row 9)
$ws.Range("A9").Value = "ModeratorPanelRestaurantsPage"
$ws.Range("B9").Value = "/moderator/restaurants"
$ws.Range("C9").Value = "Anton Tsvihun"

# Update the current selection to match the new active cell
$ws.Range("C9").Select()
